# Apply the commit's changes to GossA-HW50.xlsx:
# 1. Rename the sheet from "GossA-HW50.xpc" to "GossA".
# 2. Add a new row 16 (index 14, "HexGrid-60degTilt5degRes") mirroring the
#    existing row 15, with columns C:P all equal to 1.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Rename the worksheet
$ws.Name = "GossA"

# Duplicate row 15's formatting/content into row 16, then fix up the index
# in column A so the new row reads 14 (matching the 0-based scheme row).
$ws.Range("A15:P15").Copy($ws.Range("A16:P16"))
$ws.Cells.Item(16, 1).Value = 14
